# Work Profile and new tenant support
# Appends new OCR registration-history rows to the "AMSIN" and "AMS"
# worksheets, and corrects the run-time value on AMS row 31.

$wb = $excel.ActiveWorkbook

function Add-HistoryRow {
    param($ws, $row, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken)

    # Column A holds a date-look-alike string; force text formatting first
    # so Excel does not silently coerce it into a date serial number.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $runDate
    $ws.Cells.Item($row, 1).NumberFormat = "General"

    # Column B holds the run timestamp as a date/time serial value.
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $runTime

    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken
}

# ---------------------------------------------------------------------
# AMSIN sheet: append rows 43-48
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-HistoryRow $wsAmsin 43 "2023-03-09" 44994.56446957176 "ocr174fstcycle" 42 41 1 1.36
Add-HistoryRow $wsAmsin 44 "2023-03-13" 44998.44921864583 "174ocrflow"     42 41 1 1.44
Add-HistoryRow $wsAmsin 45 "2023-03-30" 45015.70257989583 "175ocrsc"       42 41 1 1.32
Add-HistoryRow $wsAmsin 46 "2023-03-31" 45016.49385042824 "175fnlocr"      42 41 1 1.25
Add-HistoryRow $wsAmsin 47 "2023-04-06" 45022.66083915509 "176newocr"      42 41 1 3.76
Add-HistoryRow $wsAmsin 48 "2023-04-07" 45023.66861946477 "176fstocr"      41 41 0 1.28

# ---------------------------------------------------------------------
# AMS sheet: fix run-time precision on row 31, append rows 32-35
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(31, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Cells.Item(31, 2).Value = 44978.49028510417

Add-HistoryRow $wsAms 32 "2023-03-13" 44998.55753913194 "174betaocr" 42 42 0 1.5
Add-HistoryRow $wsAms 33 "2023-03-13" 44998.85900628472 "174liveocr" 42 42 0 1.37
Add-HistoryRow $wsAms 34 "2023-03-31" 45016.56741276621 "175btocr"   42 42 0 1.44
Add-HistoryRow $wsAms 35 "2023-03-31" 45016.82906138889 "175devocr"  42 41 1 1.12
